# Re-order the observation records in rows 19-26 of the "Artfynd" sheet.
#
# The underlying data rows were reshuffled (same 8 records, new row
# positions). Mapping of final row -> original row that its contents
# came from:
#   19 <- 24
#   20 <- 19
#   21 <- 25
#   22 <- 20
#   23 <- 26
#   24 <- 21
#   25 <- 22
#   26 <- 23
#
# Whole-row ranges (columns A:AY) are copied using Range.Copy(destination)
# so that cell types/formats (text vs number vs date vs boolean) are
# preserved exactly as stored, rather than read/re-typed through .Value2
# (which would coerce things like numeric-looking text "9" or date-like
# text "1976-01-01" into real numbers/dates).
#
# Because the remapping is a cyclic permutation, rows are first copied to
# a scratch staging area (well below the used range) and then copied from
# staging into their final destinations, avoiding any source row being
# overwritten before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 19
$lastRow = 26
$lastCol = "AY"
$stageOffset = 1000

# Step 1: snapshot original rows 19-26 into staging rows (1019-1026)
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRange = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $stageRow = $r + $stageOffset
    $dstCell = $ws.Range("A" + $stageRow)
    $srcRange.Copy($dstCell)
}

# Step 2: mapping of destination row -> staged (originally-sourced) row
$mapping = @{
    19 = 24
    20 = 19
    21 = 25
    22 = 20
    23 = 26
    24 = 21
    25 = 22
    26 = 23
}

foreach ($destRow in 19..26) {
    $origRow = $mapping[$destRow]
    $stageRow = $origRow + $stageOffset
    $srcRange = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $dstCell = $ws.Range("A" + $destRow)
    $srcRange.Copy($dstCell)
}

# Step 3: clear the staging rows so they don't leak into the saved sheet
$clearRange = $ws.Range("A" + (19 + $stageOffset) + ":" + $lastCol + (26 + $stageOffset))
$clearRange.Clear()
